# Update computed price columns on Sheet1:
#   I = "Precio sin iva"                (pre-tax price)
#   L = "Unnamed: 10"                   (computed/rounded price)
#   M = "Precio Final redondeo ($)"     (final rounded price)
#
# For each product row, L and M always match, and I is simply that price
# with the 21% VAT removed (rounded to 10 significant digits, i.e. 10/1.21
# -> 8.26446281). Rows without a computed price end up at 0 for all three
# columns; rows with a computed price of 10 end up with I = 8.26446281.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows whose final price becomes 0 -------------------------------------
$ws.Range("I2:I5").Value = 0
$ws.Range("L2:L5").Value = 0
$ws.Range("M2:M5").Value = 0

$ws.Range("I12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = 0

$ws.Range("I24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = 0

$ws.Range("I28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0

$ws.Range("I36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = 0

$ws.Range("I38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 0

$ws.Range("I40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = 0

$ws.Range("I48:I49").Value = 0
$ws.Range("L48:L49").Value = 0
$ws.Range("M48:M49").Value = 0

# --- Rows whose final price becomes 10 (pre-tax 10 / 1.21 = 8.26446281) ---
$ws.Range("I6:I11").Value = 8.26446281
$ws.Range("L6:L11").Value = 10
$ws.Range("M6:M11").Value = 10

$ws.Range("I13:I23").Value = 8.26446281
$ws.Range("L13:L23").Value = 10
$ws.Range("M13:M23").Value = 10

$ws.Range("I25:I27").Value = 8.26446281
$ws.Range("L25:L27").Value = 10
$ws.Range("M25:M27").Value = 10

$ws.Range("I29:I35").Value = 8.26446281
$ws.Range("L29:L35").Value = 10
$ws.Range("M29:M35").Value = 10

$ws.Range("I37").Value = 8.26446281
$ws.Range("L37").Value = 10
$ws.Range("M37").Value = 10

$ws.Range("I39").Value = 8.26446281
$ws.Range("L39").Value = 10
$ws.Range("M39").Value = 10

$ws.Range("I41:I47").Value = 8.26446281
$ws.Range("L41:L47").Value = 10
$ws.Range("M41:M47").Value = 10

$ws.Range("I50").Value = 8.26446281
$ws.Range("L50").Value = 10
$ws.Range("M50").Value = 10
